# "Natmi following Dr Hou advice": the LR-pairs sheet for Adam23-Itga5 is
# re-run including the ECs cluster as both a sending and a target cluster
# (previously only FAPs and sCs were present). This expands the 3x3 cluster
# cross-product (FAPs, sCs, ECs) x (FAPs, sCs, ECs) = 9 rows of edge
# statistics, growing the table from rows 2-7 to rows 2-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol,
# D Target cluster, E..T the various expression / specificity / edge-weight
# statistics recomputed by natmi for the 3-cluster run.
$data = @(
    @("FAPs", "Adam23", "Itga5", "FAPs", 2, 0.6666666666666666, 0.03926266666666667, 0.117788, 0.005313231574131687, 0.005313231574131686, 3, 1, 28.72417333333333, 86.17251999999999, 0.4233259107972328, 0.4233259107972328, 1.127787642862222, 10.15008878576, 0.002249228595395911, 0.002249228595395911),
    @("FAPs", "Adam23", "Itga5", "sCs", 2, 0.6666666666666666, 0.03926266666666667, 0.117788, 0.005313231574131687, 0.005313231574131686, 3, 1, 30.56986233333333, 91.709587, 0.4505269713084062, 0.4505269713084062, 1.200254314839556, 10.802288833556, 0.002393754128953745, 0.002393754128953744),
    @("FAPs", "Adam23", "Itga5", "ECs", 2, 0.6666666666666666, 0.03926266666666667, 0.117788, 0.005313231574131687, 0.005313231574131686, 3, 1, 8.559531999999999, 25.678596, 0.126147117894361, 0.126147117894361, 0.3360700517386667, 3.024630465648, 0.0006702488497820311, 0.000670248849782031),
    @("sCs", "Adam23", "Itga5", "FAPs", 3, 1, 4.402094666666667, 13.206284, 0.5957147173375057, 0.5957147173375056, 3, 1, 28.72417333333333, 86.17251999999999, 0.4233259107972328, 0.4233259107972328, 126.4465302350755, 1138.01877211568, 0.2521814752922157, 0.2521814752922156),
    @("sCs", "Adam23", "Itga5", "sCs", 3, 1, 4.402094666666667, 13.206284, 0.5957147173375057, 0.5957147173375056, 3, 1, 30.56986233333333, 91.709587, 0.4505269713084062, 0.4505269713084062, 134.5714279383009, 1211.142851444708, 0.2683855473659097, 0.2683855473659097),
    @("sCs", "Adam23", "Itga5", "ECs", 3, 1, 4.402094666666667, 13.206284, 0.5957147173375057, 0.5957147173375056, 3, 1, 8.559531999999999, 25.678596, 0.126147117894361, 0.126147117894361, 37.67987016636266, 339.118831497264, 0.07514769467938026, 0.07514769467938025),
    @("ECs", "Adam23", "Itga5", "FAPs", 3, 1, 2.948244666666667, 8.844734000000001, 0.3989720510883627, 0.3989720510883626, 3, 1, 28.72417333333333, 86.17251999999999, 0.4233259107972328, 0.4233259107972328, 84.68589083440889, 762.17301750968, 0.1688952069096212, 0.1688952069096212),
    @("ECs", "Adam23", "Itga5", "sCs", 3, 1, 2.948244666666667, 8.844734000000001, 0.3989720510883627, 0.3989720510883626, 3, 1, 30.56986233333333, 91.709587, 0.4505269713084062, 0.4505269713084062, 90.12743358498423, 811.1469022648581, 0.1797476698135428, 0.1797476698135427),
    @("ECs", "Adam23", "Itga5", "ECs", 3, 1, 2.948244666666667, 8.844734000000001, 0.3989720510883627, 0.3989720510883626, 3, 1, 8.559531999999999, 25.678596, 0.126147117894361, 0.126147117894361, 25.23559456816267, 227.120351113464, 0.0503291743651987, 0.05032917436519869),
)

# Write the 9 data rows starting at row 2 (row 1 is the header), one row
# per (sending cluster, target cluster) pair. This overwrites the previous
# 2x2-derived rows 2-7 in place and adds the new rows 8-10, which extends
# the sheet's used range / dimension from A1:T7 to A1:T10 automatically.
$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
